# adding data to 2.13
# Appends 18 new daily log rows (A29:T46) to Sheet1, continuing the
# existing "iPhone Screen Time" tracking table (previously A1:T28).
# Column layout per header row 1:
#   A Date | B Total.ST | C Total.ST.min (formula) | D Social.ST |
#   E Social.ST.min (formula) | F Pickups | G Pickup.1st | H..T misc survey metrics
# C/E replicate the existing "parse Xh Ymin text into minutes" formula
# that is already used (as a shared formula) for rows 2-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CFormula($row) {
    return "=IF(ISERROR(FIND(""h"", B$row)), 0, LEFT(B$row, FIND(""h"", B$row)-1)*60) + IF(ISERROR(FIND(""min"", B$row)), 0, MID(B$row, IF(ISERROR(FIND(""h"", B$row)), 1, FIND(""h"", B$row)+1), FIND(""min"", B$row) - IF(ISERROR(FIND(""h"", B$row)), 1, FIND(""h"", B$row)+1)))"
}
function Get-EFormula($row) {
    return "=IF(ISERROR(FIND(""h"", D$row)), 0, LEFT(D$row, FIND(""h"", D$row)-1)*60) + IF(ISERROR(FIND(""min"", D$row)), 0, MID(D$row, IF(ISERROR(FIND(""h"", D$row)), 1, FIND(""h"", D$row)+1), FIND(""min"", D$row) - IF(ISERROR(FIND(""h"", D$row)), 1, FIND(""h"", D$row)+1)))"
}

# Each entry: row, A(date serial), B(Total.ST text), D(Social.ST text),
# F(Pickups), G(Pickup.1st time fraction), H..T(trailing numeric columns)
$rowsData = @(
    ,@(29, 45318, "59min", "22min", 34, 0.05486111111111111, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(30, 45319, "26min", "15min", 31, 0.0020833333333333333, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(31, 45320, "2h12min", "23min", 76, 0.32222222222222224, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(32, 45321, "1h34min", "37min", 85, 0.32013888888888892, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(33, 45322, "53min", "38min", 57, 0.021527777777777781, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(34, 45323, "1h10min", "24min", 84, 0.0097222222222222224, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(35, 45324, "43min", "23min", 43, 0.3125, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(36, 45325, "33min", "14min", 20, 0.34097222222222223, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(37, 45326, "6min", "1min", 23, 0.30763888888888891, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(38, 45327, "26min", "14min", 59, 0.3125, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(39, 45328, "1h29min", "1h16min", 63, 0.0062499999999999995, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(40, 45329, "33min", "20min", 70, 0, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(41, 45330, "2h", "1h27min", 78, 0.31875000000000003, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(42, 45331, "1h7min", "41min", 87, 0.20416666666666669, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(43, 45332, "1h14min", "56min", 34, 0.30416666666666664, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(44, 45333, "50min", "27min", 34, 0.3215277777777778, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(45, 45334, "1h31min", "1h2min", 66, 0.3125, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
    ,@(46, 45335, "2h3min", "1h8min", 43, 0.32916666666666666, @(0, 0, 0, 0, 1, 24, 14.5, 0, 0, 0, 6, 5, 25))
)

foreach ($r in $rowsData) {
    $row = $r[0]

    # A: date (reuse A2's date-number-format style)
    $ws.Range("A2").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $r[1]

    # B: Total.ST text
    $ws.Range("B$row").Value = $r[2]

    # C: Total.ST.min formula
    $ws.Range("C$row").Formula = Get-CFormula $row

    # D: Social.ST text
    $ws.Range("D$row").Value = $r[3]

    # E: Social.ST.min formula
    $ws.Range("E$row").Formula = Get-EFormula $row

    # F: Pickups
    $ws.Range("F$row").Value = $r[4]

    # G: Pickup.1st (reuse G2's time-number-format style)
    $ws.Range("G2").Copy($ws.Range("G$row"))
    $ws.Range("G$row").Value = $r[5]

    # H..T: trailing numeric survey columns
    $tail = $r[6]
    $tailCols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $tailCols.Length; $i++) {
        $ws.Range("$($tailCols[$i])$row").Value = $tail[$i]
    }
}

# Update sheet view: zoom to 85% and select K50 (matches workbook's saved view state)
$excel.ActiveWindow.Zoom = 85
$ws.Range("K50").Select()
